# WriteLetter: rename the "Day" column to "Id" and change its values from
# the numeric day index (1) to the literal text "Test" for the data rows.
#
# Row 11 (the "Test"/test1/test2/test3 example row) already used the text
# value "Test" in column D and numeric 2 in column A, and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Day" -> "Id"
$ws.Range("A1").Value = "Id"

# Data rows 2-10: numeric day value (1) -> text "Test"
$ws.Range("A2:A10").Value = "Test"

# Update the active selection to match the saved view state.
$ws.Range("I34").Select()
